# Updated cryptos list refresh.
# For every coin row, refresh the Price (column D) and Volume(1h) (column E)
# text values to the latest scraped snapshot. Row 44/45 also swap places
# (EnergySwap now ranks above Decentraland) with their Coin name / Link /
# Price / Volume updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds plain text such as "27.489.58" or "1.010" that must
# NOT be reinterpreted as a number (that would silently eat the grouping dot
# or a trailing zero). Force Text format on each Price cell before writing it.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextValue $ws.Range("D2") "27.489.58"
$ws.Range("E2").Value = "  +1.82%  "
Set-TextValue $ws.Range("D3") "1.859.70"
$ws.Range("E3").Value = "  +0.74%  "
Set-TextValue $ws.Range("D4") "1.010"
Set-TextValue $ws.Range("D5") "310.81"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E6").Value = "  -0.32%  "
Set-TextValue $ws.Range("D7") "0.4767"
$ws.Range("E7").Value = "  +0.00%  "
Set-TextValue $ws.Range("D8") "0.3789"
$ws.Range("E8").Value = "  +2.87%  "
Set-TextValue $ws.Range("D9") "0.07321"
$ws.Range("E9").Value = "  +1.25%  "
Set-TextValue $ws.Range("D10") "0.9293"
$ws.Range("E10").Value = "  -0.29%  "
Set-TextValue $ws.Range("D11") "20.67"
$ws.Range("E11").Value = "  +4.14%  "
Set-TextValue $ws.Range("D12") "0.07784"
$ws.Range("E12").Value = "  +0.58%  "
Set-TextValue $ws.Range("D13") "1.851.59"
$ws.Range("E13").Value = "  -0.34%  "
Set-TextValue $ws.Range("D14") "5.446"
$ws.Range("E14").Value = "  +1.10%  "
Set-TextValue $ws.Range("D15") "6.559"
$ws.Range("E15").Value = "  +1.38%  "
Set-TextValue $ws.Range("D16") "90.16"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  -0.48%  "
Set-TextValue $ws.Range("D18") "0.000008820"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("E19").Value = "  -0.42%  "
Set-TextValue $ws.Range("D20") "27.490.36"
$ws.Range("E20").Value = "  +1.68%  "
Set-TextValue $ws.Range("D21") "14.63"
$ws.Range("E21").Value = "  +0.41%  "
Set-TextValue $ws.Range("D22") "5.088"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +0.43%  "
Set-TextValue $ws.Range("D24") "1.937"
$ws.Range("E24").Value = "  -1.24%  "
Set-TextValue $ws.Range("D25") "155.22"
$ws.Range("E25").Value = "  +1.38%  "
Set-TextValue $ws.Range("D26") "18.47"
$ws.Range("E26").Value = "  +1.13%  "
Set-TextValue $ws.Range("D27") "2.005"
$ws.Range("E27").Value = "  -0.24%  "
Set-TextValue $ws.Range("D28") "115.29"
Set-TextValue $ws.Range("D29") "4.947"
$ws.Range("E29").Value = "  -0.49%  "
Set-TextValue $ws.Range("D30") "0.08867"
$ws.Range("E30").Value = "  -0.03%  "
Set-TextValue $ws.Range("D31") "3.329"
$ws.Range("E31").Value = "  +0.08%  "
Set-TextValue $ws.Range("D32") "1.202"
$ws.Range("E32").Value = "  +1.98%  "
Set-TextValue $ws.Range("D33") "0.7515"
$ws.Range("E33").Value = "  +1.03%  "
$ws.Range("E34").Value = "  +1.59%  "
Set-TextValue $ws.Range("D35") "2.705"
$ws.Range("E35").Value = "  +0.56%  "
Set-TextValue $ws.Range("D36") "0.02044"
$ws.Range("E36").Value = "  +4.08%  "
Set-TextValue $ws.Range("D37") "1.121"
$ws.Range("E37").Value = "  +0.50%  "
Set-TextValue $ws.Range("D38") "0.5546"
$ws.Range("E38").Value = "  +5.59%  "
Set-TextValue $ws.Range("D39") "0.05294"
$ws.Range("E39").Value = "  +0.52%  "
Set-TextValue $ws.Range("D40") "2.979"
$ws.Range("E40").Value = "  +0.22%  "
Set-TextValue $ws.Range("D41") "7.030"
$ws.Range("E41").Value = "  +0.03%  "
Set-TextValue $ws.Range("D42") "8.552"
$ws.Range("E42").Value = "  +2.86%  "
Set-TextValue $ws.Range("D43") "0.1519"
$ws.Range("E43").Value = "  +0.44%  "

# Row 44 / 45 swap: EnergySwap moves up to rank 44, Decentraland drops to 45.
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D44") "10.71"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.4863"
$ws.Range("E45").Value = "  +2.68%  "

Set-TextValue $ws.Range("D46") "1.010"
$ws.Range("E46").Value = "  -0.40%  "
Set-TextValue $ws.Range("D47") "103.92"
$ws.Range("E47").Value = "  +1.95%  "
$ws.Range("E48").Value = "  +3.39%  "
Set-TextValue $ws.Range("D49") "67.30"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("E51").Value = "  +1.94%  "
